$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 417.22223
$ws.Range("I6").Value = 36.57143
$ws.Range("J6").Value = 1749.5
$ws.Range("K6").Value = 109.71429
$ws.Range("L6").Value = 5248.5
$ws.Range("M6").Value = 2.285709999999995
$ws.Range("N6").Value = -5472.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 936.4
$ws.Range("I8").Value = 70
$ws.Range("J8").Value = 1802.8
$ws.Range("K8").Value = 210
$ws.Range("L8").Value = 5408.4
$ws.Range("M8").Value = -71
$ws.Range("N8").Value = -5686.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 10
$ws.Range("I31").Value = 10
$ws.Range("K31").Value = 30
$ws.Range("M31").Value = 200

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 796
$ws.Range("J32").Value = 796
$ws.Range("L32").Value = 796
$ws.Range("N32").Value = -1448

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2373.2856
$ws.Range("J86").Value = 2560
$ws.Range("L86").Value = 2560
$ws.Range("N86").Value = -4806

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1108.5
$ws.Range("I88").Value = 871
$ws.Range("K88").Value = 871
$ws.Range("M88").Value = -465

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2373.2856
$ws.Range("J89").Value = 2560
$ws.Range("L89").Value = 12800
$ws.Range("N89").Value = -24032

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1108.5
$ws.Range("I91").Value = 871
$ws.Range("K91").Value = 871
$ws.Range("M91").Value = 533

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1500
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8535.5
$ws.Range("I113").Value = 8438.125
$ws.Range("J113").Value = 8925
$ws.Range("K113").Value = 8438.125
$ws.Range("L113").Value = 8925
$ws.Range("M113").Value = -5184.125
$ws.Range("N113").Value = -15433

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3000
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("K10").Value = 2000
$ws.Range("M10").Value = -1830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1371.375
$ws.Range("I32").Value = 1186.8334
$ws.Range("K32").Value = 1186.8334
$ws.Range("M32").Value = -899.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2941.4707
$ws.Range("I45").Value = 2161.182
$ws.Range("K45").Value = 2161.182
$ws.Range("M45").Value = -1784.182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2823.8
$ws.Range("I61").Value = 2424.75
$ws.Range("K61").Value = 2424.75
$ws.Range("M61").Value = -2212.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4817.1665
$ws.Range("I74").Value = 4418.2
$ws.Range("K74").Value = 4418.2
$ws.Range("M74").Value = -3544.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4817.1665
$ws.Range("I77").Value = 4418.2
$ws.Range("K77").Value = 22091
$ws.Range("M77").Value = -17723

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2823.8
$ws.Range("I136").Value = 2424.75
$ws.Range("K136").Value = 7274.25
$ws.Range("M136").Value = -4724.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 666.2
$ws.Range("I94").Value = 666.2
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 666.2
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -215.2
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 200001310
$ws.Range("I99").Value = 250001140
$ws.Range("K99").Value = 250001140
$ws.Range("M99").Value = -249999642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1732.15
$ws.Range("I134").Value = 1315.7222
$ws.Range("J134").Value = 5480
$ws.Range("K134").Value = 3947.1666
$ws.Range("L134").Value = 16440
$ws.Range("M134").Value = -1412.1666
$ws.Range("N134").Value = -21510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1155.2
$ws.Range("I5").Value = 273.75
$ws.Range("J5").Value = 1742.8334
$ws.Range("K5").Value = 273.75
$ws.Range("L5").Value = 1742.8334
$ws.Range("M5").Value = -161.75
$ws.Range("N5").Value = -1966.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1226.4546
$ws.Range("I22").Value = 671.8889
$ws.Range("K22").Value = 671.8889
$ws.Range("M22").Value = -321.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4755.5
$ws.Range("I31").Value = 1509.4286
$ws.Range("K31").Value = 1509.4286
$ws.Range("M31").Value = -1214.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4755.5
$ws.Range("I34").Value = 1509.4286
$ws.Range("K34").Value = 1509.4286
$ws.Range("M34").Value = -1307.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 16867.334
$ws.Range("J9").Value = 300
$ws.Range("L9").Value = 900
$ws.Range("N9").Value = -1348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 4454.5454
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 6000
$ws.Range("M18").Value = -5831

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 5232.2
$ws.Range("I87").Value = 653.6667
$ws.Range("J87").Value = 12100
$ws.Range("K87").Value = 1961.0001
$ws.Range("L87").Value = 36300
$ws.Range("M87").Value = -713.0001
$ws.Range("N87").Value = -38796

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 5232.2
$ws.Range("I90").Value = 653.6667
$ws.Range("J90").Value = 12100
$ws.Range("K90").Value = 5883.0003
$ws.Range("L90").Value = 108900
$ws.Range("M90").Value = 356.9997000000003
$ws.Range("N90").Value = -121380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2581.8
$ws.Range("I114").Value = 5000
$ws.Range("J114").Value = 1977.25
$ws.Range("K114").Value = 15000
$ws.Range("L114").Value = 5931.75
$ws.Range("M114").Value = -11746
$ws.Range("N114").Value = -12439.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1486
$ws.Range("I121").Value = 476.66666
$ws.Range("J121").Value = 3000
$ws.Range("K121").Value = 1429.99998
$ws.Range("L121").Value = 9000
$ws.Range("M121").Value = -119.9999800000001
$ws.Range("N121").Value = -11620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2495.8
$ws.Range("I131").Value = 1464.75
$ws.Range("J131").Value = 3183.1667
$ws.Range("K131").Value = 4394.25
$ws.Range("L131").Value = 9549.500100000001
$ws.Range("M131").Value = 645.75
$ws.Range("N131").Value = -19629.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 35000
$ws.Range("J26").Value = 35000
$ws.Range("L26").Value = 35000
$ws.Range("N26").Value = -35560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 35000
$ws.Range("J50").Value = 35000
$ws.Range("L50").Value = 35000
$ws.Range("N50").Value = -35996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 615
$ws.Range("I97").Value = 320
$ws.Range("K97").Value = 320
$ws.Range("M97").Value = 176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 724.5
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 724.5
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 724.5
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1004.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 26478.666
$ws.Range("I51").Value = 24718
$ws.Range("J51").Value = 30000
$ws.Range("K51").Value = 24718
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = -24208
$ws.Range("N51").Value = -31020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 2983.6667
$ws.Range("I61").Value = 1560.4
$ws.Range("J61").Value = 10100
$ws.Range("K61").Value = 1560.4
$ws.Range("L61").Value = 10100
$ws.Range("M61").Value = -1268.4
$ws.Range("N61").Value = -10684

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 59999.5
$ws.Range("J115").Value = 59999.5
$ws.Range("L115").Value = 59999.5
$ws.Range("N115").Value = -63133.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2388.5
$ws.Range("I132").Value = 1777
$ws.Range("K132").Value = 5331
$ws.Range("M132").Value = -2801
